# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2..36, replacing the old Strike# values.
$kValues = @{
    2  = 7
    3  = 4
    4  = 10
    5  = 2
    6  = 2
    7  = 4
    8  = 4
    9  = 5
    10 = 3
    11 = 3
    12 = 4
    13 = 2
    14 = 6
    15 = 8
    16 = 7
    17 = 6
    18 = 5
    19 = 6
    20 = 4
    21 = 4
    22 = 2
    23 = 5
    24 = 3
    25 = 6
    26 = 7
    27 = 9
    28 = 3
    29 = 3
    30 = 6
    31 = 4
    32 = 3
    33 = 5
    34 = 4
    35 = 1
    36 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
